$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "hero defends" algorithm block: mirror the existing B:C (monster-attacks /
# hero-deffence) block into new M:N (and R:S for the wrapped note) columns,
# with two tweaked strings for the new "armor damage = 5%" defence variant.
# ---------------------------------------------------------------------------

# row 13 - sub headers
$ws.Range("M13").Value = "1)"
$ws.Range("B13").Copy()
$ws.Range("M13").PasteSpecial(-4122)

$ws.Range("N13").Value = "Monster attacks"
$ws.Range("C13").Copy()
$ws.Range("N13").PasteSpecial(-4122)

# row 15 - big explanation row
$ws.Range("M15").Value = "monster potential dmg"
$ws.Range("N15").Value = "moster.strength"

$ws.Range("R15").Value = "missDiceRoll = random number 0 to 9 "
$ws.Range("G15").Copy()
$ws.Range("R15").PasteSpecial(-4122)

$ws.Range("S15").Value = "// if (0 to 3) {attack missed} `n// if ( 4 to 5) {attack potential 50 to 70 %} `n// if ( 6 to 7) {attack potential 71 to 90 %} `n// if ( 8 to 9) { attack potential 91 to 100 % } "
$ws.Range("H15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

# row 17
$ws.Range("M17").Value = "actual monster attack"
$ws.Range("B17").Copy()
$ws.Range("M17").PasteSpecial(-4122)

$ws.Range("N17").Value = "moster.strength * missDiceRoll"

# row 19
$ws.Range("M19").Value = "2)"
$ws.Range("B19").Copy()
$ws.Range("M19").PasteSpecial(-4122)

$ws.Range("N19").Value = "hero receives damage"
$ws.Range("C19").Copy()
$ws.Range("N19").PasteSpecial(-4122)

# row 21
$ws.Range("M21").Value = "hero deffence"

# row 23/24/26 - the new "5%" armor-damage defence variant (keep this order
# so "armor damage = 5%" lands in the shared-string table before the longer
# hero.defence formula, matching the source edit)
$ws.Range("M23").Value = "if `"hero deffence`" < `"actual monster attack`""

$ws.Range("N24").Value = "hero received dmg = actual moster attack - hero deffence"

$ws.Range("N26").Value = "armor damage = 5%"

$ws.Range("N21").Value = "hero.defence + ( hero.strength + 50%) + hero.armor.head + hero.armor.shield + hero.armor.chestPlate + hero.armor.shoes"

# row 29
$ws.Range("M29").Value = "IF hero health is 0 or lower"
$ws.Range("B29").Copy()
$ws.Range("M29").PasteSpecial(-4122)

$ws.Range("N29").Value = "HERO DIES, GAME ENDS (show game results)"
$ws.Range("C29").Copy()
$ws.Range("N29").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# column widths for the new columns
# ---------------------------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 22
$ws.Columns.Item(19).ColumnWidth = 20.333333333333332

# ---------------------------------------------------------------------------
# view / selection
# ---------------------------------------------------------------------------
$ws.Range("J23").Select()
